# Update parameter files with WEP scaling and WFF_or_Ben
#
# Adds two new parameter rows to the "Parameters" sheet:
#   56: MFTC_WEP_scaling  | 1   | How should the Winter Energy Payment be scaled? ...
#   57: WFF_or_Benefit     | Max | What work decision should we assume? ...
#
# The new rows reuse the plain (no fill / no border, left-aligned) style that
# is already present on the sheet (visible on the empty header cells C1/D1),
# rather than the banded section styling used by the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write the new parameter values -----------------------------------
# All "Value"-column entries on this sheet are stored as text (e.g. "2014",
# "0.017"), not numbers, so force text format before writing "1" - otherwise
# it would be auto-detected as a number.
$ws.Range("C56").NumberFormat = "@"
$ws.Range("B56").Value = "MFTC_WEP_scaling"
$ws.Range("C56").Value = "1"
$ws.Range("D56").Value = "How should the Winter Energy Payment be scaled? Average week = 1, Winter week = 12/5, Summer week = 0"

$ws.Range("B57").Value = "WFF_or_Benefit"
$ws.Range("C57").Value = "Max"
$ws.Range("D57").Value = 'What work decision should we assume? Go off-benefit and receive IWTC = "WFF", stay on-benefit = "Benefit", or whichever gives a higher net income = "Max"'

# --- Match formatting to the rest of the workbook ----------------------
# C1:D1 already carry the plain "no fill / no border / left aligned" style
# that the two new rows use, so copy it across (this also overwrites the
# temporary "@" number format applied above) rather than re-deriving the
# style from font/fill/border properties by hand.
$ws.Range("C1").Copy()
$ws.Range("B56:D57").PasteSpecial(-4122)
$excel.CutCopyMode = 0
